# ==========================================================================
# Generate Report for Handoff
#
# Updates the localization-status workbook:
#  - refreshes the "current" handoff file name / timestamps (the source
#    markdown file 57ec59dd-...md is replaced by 7f832007-...md, and a new
#    handoff package hash f9e25b6e... replaces b52f2727...)
#  - appends two new dependency rows (png assets) that were handed off
#    together with the markdown file, on all three sheets (Overview,
#    zh-cn, de-de)
# ==========================================================================

$wb = $excel.ActiveWorkbook

# ---- constants -----------------------------------------------------------
$hyperlinkColor = 15570276   # RGB(0x64,0x95,0xED) -> matches the workbook's custom HyperLink font color
$dateFmt        = "yyyy-mm-dd HH:mm:ss"

$mdNew      = "7f832007-a04e-4392-a905-18b5927b5c8f.md"
$zhcnXlf    = "7f832007-a04e-4392-a905-18b5927b5c8f.f9e25b6eb3ef21ac89dc4bbdf6cf5f241d9849cc.zh-cn.xlf"
$dedeXlf    = "7f832007-a04e-4392-a905-18b5927b5c8f.f9e25b6eb3ef21ac89dc4bbdf6cf5f241d9849cc.de-de.xlf"

$png1       = "8a1ef8e1-7579-41e8-8e01-9a73ba333d20.png"
$png2       = "96860643-8c2d-40f5-ac74-e264fb77adeb.png"
$png1Target = "941969f9540ca50575e6870f1ef234b93b98d84a.png"
$png2Target = "fab89df9f75b3ff720a9b25caa41757930f57977.png"

$overviewDate = "2016-39-14 04:39:30"
$zhcnDate     = "2016-03-14 04:39:27"
$dedeDate     = "2016-03-14 04:39:30"
$zeroDate     = "0001-01-01 00:00:00"
$dependencyFrom = "e2e\$mdNew"

$srcRepoBase   = "https://github.com/OpenLocalizationTest/oltest/blob/6578887ce26640f75b2e7d9fcfcea7825b8ce783/e2e"
$zhcnRepoBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a6cbfc3ecc3de64ff5002aa501ff23de03c0cf02/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$dedeRepoBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/840c6c3b2a76adb94a043734faf67bb58383466d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

function Clear-AllHyperlinks($ws) {
    if ($ws.Hyperlinks.Count -gt 0) {
        $ws.Hyperlinks.Item(1).Range.Hyperlinks.Delete()
    }
}

function Add-HyperlinkCell($ws, $row, $col, $text, $url) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value2 = $text
    $ws.Hyperlinks.Add($c, $url, [System.Type]::Missing, [System.Type]::Missing, $text) | Out-Null
    $c.Font.Underline = $true
    $c.Font.Color = $hyperlinkColor
}

# ===========================================================================
# Sheet "Overview"
# ===========================================================================
$ws1 = $wb.Worksheets.Item("Overview")
Clear-AllHyperlinks $ws1

# Update the existing handoff row with the new file name / timestamp
$ws1.Cells.Item(2, 4).Value2 = $overviewDate

# New dependency rows
$ws1.Cells.Item(3, 2).Value2 = "Ready for handoff"
$ws1.Cells.Item(3, 3).Value2 = "Ready for handoff"
$ws1.Cells.Item(3, 4).Value2 = $overviewDate

$ws1.Cells.Item(4, 2).Value2 = "Ready for handoff"
$ws1.Cells.Item(4, 3).Value2 = "Ready for handoff"
$ws1.Cells.Item(4, 4).Value2 = $overviewDate

Add-HyperlinkCell $ws1 2 1 $mdNew "$srcRepoBase/$mdNew"
Add-HyperlinkCell $ws1 3 1 $png1 "$srcRepoBase/$png1"
Add-HyperlinkCell $ws1 4 1 $png2 "$srcRepoBase/$png2"

# ===========================================================================
# Sheet "zh-cn"
# ===========================================================================
$ws2 = $wb.Worksheets.Item("zh-cn")
Clear-AllHyperlinks $ws2

$ws2.Cells.Item(2, 5).Value2 = $zhcnDate
$ws2.Cells.Item(2, 5).NumberFormat = $dateFmt

# Row 3 - first png dependency
$ws2.Cells.Item(3, 3).Value2 = "Ready for handoff"
$ws2.Cells.Item(3, 5).Value2 = $zhcnDate
$ws2.Cells.Item(3, 5).NumberFormat = $dateFmt
$ws2.Cells.Item(3, 8).Value2 = $zeroDate
$ws2.Cells.Item(3, 9).Value2 = "IsDependency"
$ws2.Cells.Item(3, 10).Value2 = $dependencyFrom

# Row 4 - second png dependency
$ws2.Cells.Item(4, 3).Value2 = "Ready for handoff"
$ws2.Cells.Item(4, 5).Value2 = $zhcnDate
$ws2.Cells.Item(4, 5).NumberFormat = $dateFmt
$ws2.Cells.Item(4, 8).Value2 = $zeroDate
$ws2.Cells.Item(4, 9).Value2 = "IsDependency"
$ws2.Cells.Item(4, 10).Value2 = $dependencyFrom

# Hyperlinks, created in the same left-to-right / top-to-bottom order as
# the original workbook so relationship ids line up (rId2 .. rId9)
Add-HyperlinkCell $ws2 2 1 $mdNew "$srcRepoBase/$mdNew"
Add-HyperlinkCell $ws2 2 2 ".md" "$srcRepoBase/$mdNew"
Add-HyperlinkCell $ws2 2 4 $zhcnXlf "$zhcnRepoBase/$zhcnXlf"
Add-HyperlinkCell $ws2 3 1 $png1 "$srcRepoBase/$png1"
Add-HyperlinkCell $ws2 3 2 ".png" "$srcRepoBase/$png1"
Add-HyperlinkCell $ws2 3 4 $png1Target "$zhcnRepoBase/$png1Target"
Add-HyperlinkCell $ws2 4 1 $png2 "$srcRepoBase/$png2"
Add-HyperlinkCell $ws2 4 2 ".png" "$srcRepoBase/$png2"
Add-HyperlinkCell $ws2 4 4 $png2Target "$zhcnRepoBase/$png2Target"

# ===========================================================================
# Sheet "de-de"
# ===========================================================================
$ws3 = $wb.Worksheets.Item("de-de")
Clear-AllHyperlinks $ws3

$ws3.Cells.Item(2, 5).Value2 = $dedeDate
$ws3.Cells.Item(2, 5).NumberFormat = $dateFmt

# Row 3 - first png dependency
$ws3.Cells.Item(3, 3).Value2 = "Ready for handoff"
$ws3.Cells.Item(3, 5).Value2 = $dedeDate
$ws3.Cells.Item(3, 5).NumberFormat = $dateFmt
$ws3.Cells.Item(3, 8).Value2 = $zeroDate
$ws3.Cells.Item(3, 9).Value2 = "IsDependency"
$ws3.Cells.Item(3, 10).Value2 = $dependencyFrom

# Row 4 - second png dependency
$ws3.Cells.Item(4, 3).Value2 = "Ready for handoff"
$ws3.Cells.Item(4, 5).Value2 = $dedeDate
$ws3.Cells.Item(4, 5).NumberFormat = $dateFmt
$ws3.Cells.Item(4, 8).Value2 = $zeroDate
$ws3.Cells.Item(4, 9).Value2 = "IsDependency"
$ws3.Cells.Item(4, 10).Value2 = $dependencyFrom

Add-HyperlinkCell $ws3 2 1 $mdNew "$srcRepoBase/$mdNew"
Add-HyperlinkCell $ws3 2 2 ".md" "$srcRepoBase/$mdNew"
Add-HyperlinkCell $ws3 2 4 $dedeXlf "$dedeRepoBase/$dedeXlf"
Add-HyperlinkCell $ws3 3 1 $png1 "$srcRepoBase/$png1"
Add-HyperlinkCell $ws3 3 2 ".png" "$srcRepoBase/$png1"
Add-HyperlinkCell $ws3 3 4 $png1Target "$dedeRepoBase/$png1Target"
Add-HyperlinkCell $ws3 4 1 $png2 "$srcRepoBase/$png2"
Add-HyperlinkCell $ws3 4 2 ".png" "$srcRepoBase/$png2"
Add-HyperlinkCell $ws3 4 4 $png2Target "$dedeRepoBase/$png2Target"

Write-Host "Report generated for handoff."
